$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K24").Value = -0.2401694470146267
$ws.Range("J25").Value = -0.3012619146487851
$ws.Range("K25").Value = 0.9653040704216436
$ws.Range("I26").Value = -0.4453260117243474
$ws.Range("J26").Value = 0.8212399733460813
$ws.Range("K26").Value = 0.1149727305944879
$ws.Range("H27").Value = -0.3960164397552333
$ws.Range("I27").Value = 0.8705495453151955
$ws.Range("J27").Value = 0.164282302563602
$ws.Range("K27").Value = -1.445378088853165
$ws.Range("G28").Value = -0.4646074970197469
$ws.Range("H28").Value = 0.8019584880506818
$ws.Range("I28").Value = 0.09569124529908837
$ws.Range("J28").Value = -1.513969146117679
$ws.Range("K28").Value = 2.514191877581562
$ws.Range("F29").Value = -0.3625832606467455
$ws.Range("G29").Value = 0.9039827244236832
$ws.Range("H29").Value = 0.1977154816720898
$ws.Range("I29").Value = -1.411944909744677
$ws.Range("J29").Value = 2.616216113954564
$ws.Range("K29").Value = -0.5603787180733534
$ws.Range("E30").Value = -0.3042477716424103
$ws.Range("F30").Value = 0.9623182134280184
$ws.Range("G30").Value = 0.256050970676425
$ws.Range("H30").Value = -1.353609420740342
$ws.Range("I30").Value = 2.674551602958899
$ws.Range("J30").Value = -0.5020432290690182
$ws.Range("K30").Value = 0.3875911078309982
$ws.Range("D31").Value = -0.2497932332702604
$ws.Range("E31").Value = 1.016772751800168
$ws.Range("F31").Value = 0.3105055090485749
$ws.Range("G31").Value = -1.299154882368192
$ws.Range("H31").Value = 2.729006141331049
$ws.Range("I31").Value = -0.4475886906968682
$ws.Range("J31").Value = 0.4420456462031481
$ws.Range("K31").Value = -1.833582676193962
$ws.Range("C32").Value = -0.2778442010910621
$ws.Range("D32").Value = 0.9887217839793667
$ws.Range("E32").Value = 0.2824545412277732
$ws.Range("F32").Value = -1.327205850188994
$ws.Range("G32").Value = 2.700955173510247
$ws.Range("H32").Value = -0.4756396585176699
$ws.Range("I32").Value = 0.4139946783823464
$ws.Range("J32").Value = -1.861633644014764
$ws.Range("K32").Value = 0.4824515838321392
$ws.Range("B33").Value = -0.1336561884836082
$ws.Range("C33").Value = 1.132909796586821
$ws.Range("D33").Value = 0.4266425538352271
$ws.Range("E33").Value = -1.18301783758154
$ws.Range("F33").Value = 2.845143186117701
$ws.Range("G33").Value = -0.331451645910216
$ws.Range("H33").Value = 0.5581826909898003
$ws.Range("I33").Value = -1.71744563140731
$ws.Range("J33").Value = 0.6266395964395931
$ws.Range("K33").Value = 0.003271690537692051
$ws.Range("B34").Value = 1.307066490313564
$ws.Range("C34").Value = 0.6007992475619707
$ws.Range("D34").Value = -1.008861143854797
$ws.Range("E34").Value = 3.019299879844445
$ws.Range("F34").Value = -0.1572949521834724
$ws.Range("G34").Value = 0.7323393847165439
$ws.Range("H34").Value = -1.543288937680567
$ws.Range("I34").Value = 0.8007962901663367
$ws.Range("J34").Value = 0.1774283842644356
$ws.Range("K34").Value = 5.68814648014622
$ws.Range("B35").Value = 0.4257766076739633
$ws.Range("C35").Value = -1.183883783742804
$ws.Range("D35").Value = 2.844277239956437
$ws.Range("E35").Value = -0.3323175920714798
$ws.Range("F35").Value = 0.5573167448285365
$ws.Range("G35").Value = -1.718311577568574
$ws.Range("H35").Value = 0.6257736502783293
$ws.Range("I35").Value = 0.002405744376428265
$ws.Range("J35").Value = 5.513123840258213
$ws.Range("K35").Value = 0.548303282438791
$ws.Range("B36").Value = -1.089199807152462
$ws.Range("C36").Value = 2.938961216546779
$ws.Range("D36").Value = -0.2376336154811378
$ws.Range("E36").Value = 0.6520007214188785
$ws.Range("F36").Value = -1.623627600978232
$ws.Range("G36").Value = 0.7204576268686713
$ws.Range("H36").Value = 0.09708972096677027
$ws.Range("I36").Value = 5.607807816848555
$ws.Range("J36").Value = 0.642987259029133
$ws.Range("K36").Value = 0.4972523523676307
$ws.Range("B37").Value = 2.991208970153319
$ws.Range("C37").Value = -0.1853858618745983
$ws.Range("D37").Value = 0.7042484750254181
$ws.Range("E37").Value = -1.571379847371692
$ws.Range("F37").Value = 0.7727053804752109
$ws.Range("G37").Value = 0.1493374745733098
$ws.Range("H37").Value = 5.660055570455095
$ws.Range("I37").Value = 0.6952350126356726
$ws.Range("J37").Value = 0.5495001059741702
$ws.Range("K37").Value = 0.4253868890927299
$ws.Range("B38").Value = -0.2633354489934661
$ws.Range("C38").Value = 0.6262988879065503
$ws.Range("D38").Value = -1.64932943449056
$ws.Range("E38").Value = 0.6947557933563431
$ws.Range("F38").Value = 0.07138788745444202
$ws.Range("G38").Value = 5.582105983336227
$ws.Range("H38").Value = 0.6172854255168048
$ws.Range("I38").Value = 0.4715505188553024
$ws.Range("J38").Value = 0.3474373019738621
$ws.Range("K38").Value = 1.121628485305244
$ws.Range("B39").Value = 0.6469022353042405
$ws.Range("C39").Value = -1.62872608709287
$ws.Range("D39").Value = 0.7153591407540333
$ws.Range("E39").Value = 0.09199123485213223
$ws.Range("F39").Value = 5.602709330733918
$ws.Range("G39").Value = 0.637888772914495
$ws.Range("H39").Value = 0.4921538662529926
$ws.Range("I39").Value = 0.3680406493715523
$ws.Range("J39").Value = 1.142231832702934
$ws.Range("K39").Value = -0.31932174300057
$ws.Range("B40").Value = -1.853005556311659
$ws.Range("C40").Value = 0.4910796715352439
$ws.Range("D40").Value = -0.1322882343666572
$ws.Range("E40").Value = 5.378429861515128
$ws.Range("F40").Value = 0.4136093036957056
$ws.Range("G40").Value = 0.2678743970342032
$ws.Range("H40").Value = 0.1437611801527629
$ws.Range("I40").Value = 0.9179523634841451
$ws.Range("J40").Value = -0.5436012122193594
$ws.Range("K40").Value = 0.05670714581054659
$ws.Range("B41").Value = 0.8890957203403786
$ws.Range("C41").Value = 0.2657278144384776
$ws.Range("D41").Value = 5.776445910320263
$ws.Range("E41").Value = 0.8116253525008403
$ws.Range("F41").Value = 0.665890445839338
$ws.Range("G41").Value = 0.5417772289578977
$ws.Range("H41").Value = 1.31596841228928
$ws.Range("I41").Value = -0.1455851634142247
$ws.Range("J41").Value = 0.4547231946156813
$ws.Range("K41").Value = 0.4730809959422544
$ws.Range("B42").Value = -0.1503833582008396
$ws.Range("C42").Value = 5.360334737680946
$ws.Range("D42").Value = 0.3955141798615231
$ws.Range("E42").Value = 0.2497792732000207
$ws.Range("F42").Value = 0.1256660563185804
$ws.Range("G42").Value = 0.8998572396499627
$ws.Range("H42").Value = -0.5616963360535419
$ws.Range("I42").Value = 0.0386120219763641
$ws.Range("J42").Value = 0.05696982330293715
$ws.Range("K42").Value = 0.8153243866718058
$ws.Range("B43").Value = 5.178024097175236
$ws.Range("C43").Value = 0.2132035393558141
$ws.Range("D43").Value = 0.06746863269431172
$ws.Range("E43").Value = -0.05664458418712859
$ws.Range("F43").Value = 0.7175465991442536
$ws.Range("G43").Value = -0.744006976559251
$ws.Range("H43").Value = -0.1436986185293449
$ws.Range("I43").Value = -0.1253408172027718
$ws.Range("J43").Value = 0.6330137461660968
$ws.Range("K43").Value = -0.5012717062398651
$ws.Range("B44").Value = -0.4043520997521085
$ws.Range("C44").Value = -0.5500870064136109
$ws.Range("D44").Value = -0.6742002232950512
$ws.Range("E44").Value = 0.09999096003633101
$ws.Range("F44").Value = -1.361562615667173
$ws.Range("G44").Value = -0.7612542576372675
$ws.Range("H44").Value = -0.7428964563106945
$ws.Range("I44").Value = 0.01545810705817424
$ws.Range("J44").Value = -1.118827345347788
$ws.Range("B45").Value = -0.259964220622777
$ws.Range("C45").Value = -0.3840774375042173
$ws.Range("D45").Value = 0.3901137458271649
$ws.Range("E45").Value = -1.07143982987634
$ws.Range("F45").Value = -0.4711314718464336
$ws.Range("G45").Value = -0.4527736705198606
$ws.Range("H45").Value = 0.3055808928490081
$ws.Range("I45").Value = -0.8287045595569538
$ws.Range("B46").Value = -0.2835772387253058
$ws.Range("C46").Value = 0.4906139446060763
$ws.Range("D46").Value = -0.9709396310974282
$ws.Range("E46").Value = -0.3706312730675222
$ws.Range("F46").Value = -0.3522734717409491
$ws.Range("G46").Value = 0.4060810916279196
$ws.Range("H46").Value = -0.7282043607780423
$ws.Range("B47").Value = 0.1465102285875062
$ws.Range("C47").Value = -1.315043347115998
$ws.Range("D47").Value = -0.7147349890860923
$ws.Range("E47").Value = -0.6963771877595193
$ws.Range("F47").Value = 0.06197737560934941
$ws.Range("G47").Value = -1.072308076796612
$ws.Range("B48").Value = -0.9587741681389913
$ws.Range("C48").Value = -0.3584658101090852
$ws.Range("D48").Value = -0.3401080087825122
$ws.Range("E48").Value = 0.4182465545863565
$ws.Range("F48").Value = -0.7160388978196054
$ws.Range("B49").Value = -0.2401767394741711
$ws.Range("C49").Value = -0.2218189381475981
$ws.Range("D49").Value = 0.5365356252212706
$ws.Range("E49").Value = -0.5977498271846913
$ws.Range("B50").Value = -0.219261153323231
$ws.Range("C50").Value = 0.5390934100456377
$ws.Range("D50").Value = -0.5951920423603241
$ws.Range("B51").Value = 0.7152142819569747
$ws.Range("C51").Value = -0.4190711704489871
$ws.Range("B52").Value = -0.377607739757282
